# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values with the newly curated dimension names
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("B3").Value = "dim"
$ws.Range("B4").Value = "URI-Municipio"
$ws.Range("D4").Value = "URI-Comunidad"

# Remove the old aragon mapping reference entirely (row 5, column D)
$ws.Range("D5").Clear()
